# Alt Index - Initial.xlsx : "Update files and fix bugs"
#
# 1. Column A (rows 2-21) gets refreshed unique identifiers (new random
#    hex suffixes on the existing R-20x-Ix-... codes).
# 2. A brand-new row 22 is appended, continuing the same pattern
#    (R-222-I22-CS2-7D79 / AltIndex-V / formulas / Good / ENG / Engine).
# 3. The used range grows from A1:H21 to A1:H22 (automatic once the new
#    row has data) and the active selection moves to G27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- refreshed identifiers in column A (rows 2-21) -------------------------
# Ordered (row, new-id) pairs so the writes happen top-to-bottom, same as the
# natural row order in the sheet.
$newIds = @(
    , @(2,  "R-201-I1-CS1-674B")
    , @(3,  "R-202-I2-Cb1-8C22")
    , @(4,  "R-203-I3-CS2-7A3F")
    , @(5,  "R-204-I4-Cf1-3967")
    , @(6,  "R-205-I5-Ce1-236C")
    , @(7,  "R-206-I6-CT1-012B")
    , @(8,  "R-207-I7-CE2-4D7F")
    , @(9,  "R-208-I8-CS1-22BF")
    , @(10, "R-209-I9-Cb1-6694")
    , @(11, "R-210-I10-CS2-BD7D")
    , @(12, "R-211-I11-Cf1-9B8C")
    , @(13, "R-212-I12-Ce1-2186")
    , @(14, "R-213-I13-CT1-A74F")
    , @(15, "R-214-I14-CE2-CF3C")
    , @(16, "R-215-I15-CS1-96FB")
    , @(17, "R-216-I16-Cb1-F7E7")
    , @(18, "R-217-I17-CS2-96CC")
    , @(19, "R-218-I18-Cf1-E368")
    , @(20, "R-219-I19-Ce1-922D")
    , @(21, "R-220-I20-CT1-C738")
)

foreach ($pair in $newIds) {
    $row = $pair[0]
    $id  = $pair[1]
    $ws.Cells.Item($row, 1).Value = $id
}

# --- new row 22 -------------------------------------------------------------
$ws.Cells.Item(22, 1).Value = "R-222-I22-CS2-7D79"
$ws.Cells.Item(22, 2).Value = "AltIndex-V"
$ws.Cells.Item(22, 3).Formula = "=5*ROW()"
$ws.Cells.Item(22, 4).Formula = "=100*ROW()"
$ws.Cells.Item(22, 5).Value = "Good"
$ws.Cells.Item(22, 6).Value = "ENG"
$ws.Cells.Item(22, 7).Value = "Engine"

# --- selection / window position mirror the saved view in the target file -
$ws.Range("G27").Select()

$win = $excel.Windows.Item(1)
$win.Left = 5175
$win.Top = 3570
